$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.793.89'
$ws.Range('D3').Value = '1.700.09'
$ws.Range('E3').Value = '  +0.13%  '
$c = $ws.Range('D4')
$c.Formula = '="1.002"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E4').Value = '  +0.15%  '
$c = $ws.Range('D5')
$c.Formula = '="318.66"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E5').Value = '  +1.26%  '
$c = $ws.Range('D6')
$c.Formula = '="1.001"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E6').Value = '  +0.07%  '
$c = $ws.Range('D7')
$c.Formula = '="0.3956"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E7').Value = '  +0.25%  '
$c = $ws.Range('D8')
$c.Formula = '="0.4057"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E8').Value = '  +0.24%  '
$c = $ws.Range('D9')
$c.Formula = '="1.502"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E9').Value = '  -1.92%  '
$c = $ws.Range('D10')
$c.Formula = '="1.003"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E10').Value = '  +0.22%  '
$c = $ws.Range('D11')
$c.Formula = '="52.94"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E11').Value = '  -9.37%  '
$c = $ws.Range('D12')
$c.Formula = '="0.08893"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E12').Value = '  +0.87%  '
$c = $ws.Range('D13')
$c.Formula = '="7.321"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E13').Value = '  +1.17%  '
$c = $ws.Range('D14')
$c.Formula = '="23.56"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E14').Value = '  +0.67%  '
$c = $ws.Range('D15')
$c.Formula = '="8.030"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E15').Value = '  +5.59%  '
$c = $ws.Range('D16')
$c.Formula = '="0.00001321"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E16').Value = '  -1.07%  '
$ws.Range('D17').Value = '1.700.64'
$ws.Range('E17').Value = '  +0.21%  '
$c = $ws.Range('D18')
$c.Formula = '="100.18"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E18').Value = '  -0.71%  '
$c = $ws.Range('D19')
$c.Formula = '="0.07040"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E19').Value = '  -0.56%  '
$c = $ws.Range('D20')
$c.Formula = '="19.68"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E20').Value = '  +0.12%  '
$c = $ws.Range('D21')
$c.Formula = '="7.049"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E21').Value = '  +4.30%  '
$c = $ws.Range('D22')
$c.Formula = '="1.001"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E22').Value = '  +0.00%  '
$c = $ws.Range('D23')
$c.Formula = '="14.45"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E23').Value = '  +1.69%  '
$ws.Range('D24').Value = '24.777.36'
$ws.Range('E24').Value = '  +0.51%  '
$c = $ws.Range('D25')
$c.Formula = '="3.211"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E25').Value = '  +7.26%  '
$c = $ws.Range('D26')
$c.Formula = '="2.364"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E26').Value = '  +1.93%  '
$c = $ws.Range('D27')
$c.Formula = '="22.79"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E27').Value = '  +1.42%  '
$c = $ws.Range('D28')
$c.Formula = '="162.29"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E28').Value = '  +1.23%  '
$c = $ws.Range('D29')
$c.Formula = '="137.33"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E29').Value = '  +2.28%  '
$c = $ws.Range('D30')
$c.Formula = '="5.185"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E30').Value = '  -0.26%  '
$c = $ws.Range('D31')
$c.Formula = '="7.706"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E31').Value = '  +5.08%  '
$c = $ws.Range('D32')
$c.Formula = '="0.08710"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E32').Value = '  +0.98%  '
$c = $ws.Range('D33')
$c.Formula = '="1.072"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E33').Value = '  -3.64%  '
$c = $ws.Range('D34')
$c.Formula = '="7.136"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E34').Value = '  -3.55%  '
$c = $ws.Range('D35')
$c.Formula = '="11.38"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E35').Value = '  +2.65%  '
$c = $ws.Range('D36')
$c.Formula = '="0.2748"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('B37').Value = 'WEMIXTOKEN'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range('D37')
$c.Formula = '="1.918"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E37').Value = '  -3.65%  '
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range('D38')
$c.Formula = '="14.55"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E38').Value = '  -2.11%  '
$c = $ws.Range('D39')
$c.Formula = '="0.09227"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E39').Value = '  +1.53%  '
$c = $ws.Range('D40')
$c.Formula = '="0.02731"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E40').Value = '  -1.64%  '
$c = $ws.Range('D41')
$c.Formula = '="1.468"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E41').Value = '  -0.24%  '
$c = $ws.Range('D42')
$c.Formula = '="0.7731"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E42').Value = '  -0.44%  '
$c = $ws.Range('D43')
$c.Formula = '="16.24"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E43').Value = '  +4.47%  '
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c = $ws.Range('D44')
$c.Formula = '="0.7207"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E44').Value = '  -0.91%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D45')
$c.Formula = '="2.571"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E45').Value = '  +2.99%  '
$c = $ws.Range('D46')
$c.Formula = '="4.240"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E46').Value = '  +1.07%  '
$c = $ws.Range('D47')
$c.Formula = '="1.001"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E47').Value = '  +0.07%  '
$c = $ws.Range('D48')
$c.Formula = '="140.69"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E48').Value = '  -0.44%  '
$c = $ws.Range('D49')
$c.Formula = '="1.325"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E49').Value = '  +2.52%  '
$c = $ws.Range('D50')
$c.Formula = '="91.77"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E50').Value = '  +5.04%  '
$c = $ws.Range('D51')
$c.Formula = '="0.07999"'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E51').Value = '  -0.26%  '

$excel.CutCopyMode = $false
